$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The status text "Ready for handoff" is a shared string referenced by the
# Overview sheet's E3/F3 (zh-cn / de-de status columns) as well as the
# per-locale sheets' own "Status" column (C3). All occurrences flip to the
# new status text.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file name mismatch in the "Error Detail" column
# (P) for row 3 on each localized sheet.
$wsZhCn.Range("P3").Value = "Handback file name: zjowasme.21w is different with handoff file name: f762ea34-3777-4ddc-b8f3-0225b47b10bc.0e19928eea6eb1aff0c2fb2aebd5193f5e790349.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: zjowasme.21w is different with handoff file name: f762ea34-3777-4ddc-b8f3-0225b47b10bc.0e19928eea6eb1aff0c2fb2aebd5193f5e790349.de-de."

# Widen the "Error Detail" column (P, the 16th column) on both localized
# sheets from ~13.75 to 40 characters so the new message is readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 40 - 5/6
$wsDeDe.Columns.Item(16).ColumnWidth = 40 - 5/6
